# Auto-generated edit script applying numeric corrections to H:N columns
# across multiple leve-profit worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ALC!row137 - Magnesia Whetstone
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 886.34485
$ws.Range("I137").Value = 844
$ws.Range("J137").Value = 997.5
$ws.Range("K137").Value = 2532
$ws.Range("L137").Value = 2992.5
$ws.Range("M137").Value = 18
$ws.Range("N137").Value = -8092.5

# ARM!row61 - Cobalt Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1533.6923
$ws.Range("I61").Value = 1501.375
$ws.Range("K61").Value = 1501.375
$ws.Range("M61").Value = -1289.375

# ARM!row122 - High Durium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 34684.676
$ws.Range("I122").Value = 2355.8262
$ws.Range("J122").Value = 127630.125
$ws.Range("K122").Value = 7067.4786
$ws.Range("L122").Value = 382890.375
$ws.Range("M122").Value = -4617.4786
$ws.Range("N122").Value = -387790.375

# ARM!row132 - Mountain Chromite Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 45502680
$ws.Range("I132").Value = 76925384
$ws.Range("J132").Value = 114333.555
$ws.Range("K132").Value = 230776152
$ws.Range("L132").Value = 343000.665
$ws.Range("M132").Value = -230773622
$ws.Range("N132").Value = -348060.665

# ARM!row136 - Cobalt Tungsten Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1533.6923
$ws.Range("I136").Value = 1501.375
$ws.Range("K136").Value = 4504.125
$ws.Range("M136").Value = -1954.125

# BSM!row20 - Iron Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11249.5
$ws.Range("I20").Value = 9999
$ws.Range("K20").Value = 9999
$ws.Range("M20").Value = -9752

# BSM!row52 - Mythril File
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 24881.666
$ws.Range("J52").Value = 24881.666
$ws.Range("L52").Value = 24881.666
$ws.Range("N52").Value = -25407.666

# BSM!row121 - Dwarven Mythril File
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H121").Value = 24881.666
$ws.Range("J121").Value = 24881.666
$ws.Range("L121").Value = 24881.666
$ws.Range("N121").Value = -28375.666

# BSM!row134 - Ruthenium Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6260.1353
$ws.Range("I134").Value = 2036.9667
$ws.Range("J134").Value = 24359.428
$ws.Range("K134").Value = 6110.9001
$ws.Range("L134").Value = 73078.284
$ws.Range("M134").Value = -3575.9001
$ws.Range("N134").Value = -78148.284

# CRP!row31 - Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5436578
$ws.Range("I31").Value = 6758356
$ws.Range("J31").Value = 2601.5557
$ws.Range("K31").Value = 6758356
$ws.Range("L31").Value = 2601.5557
$ws.Range("M31").Value = -6758061
$ws.Range("N31").Value = -3191.5557

# CRP!row34 - Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5436578
$ws.Range("I34").Value = 6758356
$ws.Range("J34").Value = 2601.5557
$ws.Range("K34").Value = 6758356
$ws.Range("L34").Value = 2601.5557
$ws.Range("M34").Value = -6758154
$ws.Range("N34").Value = -3005.5557

# CRP!row58 - Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1097.1111
$ws.Range("I58").Value = 1002.8333
$ws.Range("J58").Value = 1285.6666
$ws.Range("K58").Value = 1002.8333
$ws.Range("L58").Value = 1285.6666
$ws.Range("M58").Value = -799.8333
$ws.Range("N58").Value = -1691.6666

# CRP!row136 - Dark Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1097.1111
$ws.Range("I136").Value = 1002.8333
$ws.Range("J136").Value = 1285.6666
$ws.Range("K136").Value = 3008.4999
$ws.Range("L136").Value = 3856.9998
$ws.Range("M136").Value = -458.4998999999998
$ws.Range("N136").Value = -8956.9998

# CUL!row12 - Kukuru Butter
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 144.16667
$ws.Range("I12").Value = 1.5
$ws.Range("J12").Value = 215.5
$ws.Range("K12").Value = 4.5
$ws.Range("L12").Value = 646.5
$ws.Range("M12").Value = 168.5
$ws.Range("N12").Value = -992.5

# CUL!row44 - Knight's Bread
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 8240.799999999999
$ws.Range("I44").Value = 300
$ws.Range("J44").Value = 40004
$ws.Range("K44").Value = 900
$ws.Range("L44").Value = 120012
$ws.Range("M44").Value = -502
$ws.Range("N44").Value = -120808

# CUL!row64 - Baked Onion Soup
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2399.3333
$ws.Range("I64").Value = 200
$ws.Range("J64").Value = 3499
$ws.Range("K64").Value = 600
$ws.Range("L64").Value = 10497
$ws.Range("M64").Value = -330
$ws.Range("N64").Value = -11037

# CUL!row67 - Baked Onion Soup
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 2399.3333
$ws.Range("I67").Value = 200
$ws.Range("J67").Value = 3499
$ws.Range("K67").Value = 600
$ws.Range("L67").Value = 10497
$ws.Range("M67").Value = 336
$ws.Range("N67").Value = -12369

# CUL!row113 - Night Vinegar
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 994.2895
$ws.Range("I113").Value = 698.6667
$ws.Range("J113").Value = 1006.43835
$ws.Range("K113").Value = 2096.0001
$ws.Range("L113").Value = 3019.31505
$ws.Range("M113").Value = 73.9998999999998
$ws.Range("N113").Value = -7359.31505

# GSM!row132 - Lar Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 394690.25
$ws.Range("I132").Value = 59997.59
$ws.Range("J132").Value = 911942.5600000001
$ws.Range("K132").Value = 179992.77
$ws.Range("L132").Value = 2735827.68
$ws.Range("M132").Value = -177462.77
$ws.Range("N132").Value = -2740887.68

# LTW!row39 - Boarskin Himantes
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 4400
$ws.Range("J39").Value = 4400
$ws.Range("L39").Value = 4400
$ws.Range("N39").Value = -5320

# LTW!row41 - Fingerless Boarskin Gloves
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 5099
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

# LTW!row132 - Silver Lobo Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 26127.61
$ws.Range("I132").Value = 51945.3
$ws.Range("J132").Value = 1539.3334
$ws.Range("K132").Value = 155835.9
$ws.Range("L132").Value = 4618.0002
$ws.Range("M132").Value = -153305.9
$ws.Range("N132").Value = -9678.0002

# LTW!row136 - Br'aax Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6110.7427
$ws.Range("I136").Value = 8028.8423
$ws.Range("J136").Value = 3833
$ws.Range("K136").Value = 24086.5269
$ws.Range("L136").Value = 11499
$ws.Range("M136").Value = -21536.5269
$ws.Range("N136").Value = -16599

# WVR!row70 - Holy Rainbow Shoes
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 11000
$ws.Range("J70").Value = 11000
$ws.Range("L70").Value = 11000
$ws.Range("N70").Value = -11630

# WVR!row73 - Holy Rainbow Shoes
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 11000
$ws.Range("J73").Value = 11000
$ws.Range("L73").Value = 11000
$ws.Range("N73").Value = -13184

# WVR!row75 - Ramie Turban of Crafting
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 14800
$ws.Range("J75").Value = 14800
$ws.Range("L75").Value = 14800
$ws.Range("N75").Value = -16672

# WVR!row78 - Ramie Turban of Crafting
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 14800
$ws.Range("J78").Value = 14800
$ws.Range("L78").Value = 44400
$ws.Range("N78").Value = -53760

# WVR!row136 - Sarcenet Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 28982.389
$ws.Range("I136").Value = 50550.3
$ws.Range("J136").Value = 2022.5
$ws.Range("K136").Value = 151650.9
$ws.Range("L136").Value = 6067.5
$ws.Range("M136").Value = -149100.9
$ws.Range("N136").Value = -11167.5
